$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell K1 "Prot", matching the style of existing header cells (e.g. J1):
# bold font, centered/top aligned, thin border all around.
$ws.Range("K1").Value = "Prot"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4160
$ws.Range("K1").Borders.LineStyle = 1

# Fill new column K (Prot) with "PEPTIDE" for all 12 data rows
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 11).Value = "PEPTIDE"
}
